$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 of Sheet1 stores its values as text (inlineStr in the source XML),
# even though several of them look numeric (e.g. "8761", "5746").
# Writing a plain numeric-looking string via COM gets auto-coerced into a
# real number cell, which would change the cell's stored type. To keep
# these as text cells (matching the rest of the row, e.g. D2/G2/H2 which
# are text), force text entry with a leading apostrophe, then restore the
# "Normal" style so no stray number-format/style is left behind on the
# cell.

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "C2" "0"
Set-TextValue "E2" "0.439"
Set-TextValue "F2" "0.007"
Set-TextValue "I2" "5726"
Set-TextValue "P2" "0.0"
Set-TextValue "Q2" "0.0"
